$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1092.7778
$ws.Range("I19").Value = 815.8333
$ws.Range("J19").Value = 1314.3334
$ws.Range("K19").Value = 815.8333
$ws.Range("L19").Value = 1314.3334
$ws.Range("M19").Value = -640.8333
$ws.Range("N19").Value = -1664.3334

$ws.Range("H62").Value = 41669236
$ws.Range("I62").Value = 55557310
$ws.Range("J62").Value = 5017.6665
$ws.Range("K62").Value = 55557310
$ws.Range("L62").Value = 5017.6665
$ws.Range("M62").Value = -55556686
$ws.Range("N62").Value = -6265.6665

$ws.Range("H65").Value = 41669236
$ws.Range("I65").Value = 55557310
$ws.Range("J65").Value = 5017.6665
$ws.Range("K65").Value = 277786550
$ws.Range("L65").Value = 25088.3325
$ws.Range("M65").Value = -277783430
$ws.Range("N65").Value = -31328.3325

$ws.Range("H129").Value = 991.2632
$ws.Range("I129").Value = 923.5
$ws.Range("J129").Value = 999.2353000000001
$ws.Range("K129").Value = 2770.5
$ws.Range("L129").Value = 2997.7059
$ws.Range("M129").Value = 2229.5
$ws.Range("N129").Value = -12997.7059

$ws.Range("H137").Value = 1091133.5
$ws.Range("I137").Value = 1451.6666
$ws.Range("J137").Value = 2647821.8
$ws.Range("K137").Value = 4354.9998
$ws.Range("L137").Value = 7943465.399999999
$ws.Range("M137").Value = -1804.9998
$ws.Range("N137").Value = -7948565.399999999

$ws.Range("H141").Value = 1001463.3
$ws.Range("I141").Value = 1035651.75
$ws.Range("J141").Value = 10000
$ws.Range("K141").Value = 3106955.25
$ws.Range("L141").Value = 30000
$ws.Range("M141").Value = -3101775.25
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1397.6052
$ws.Range("I2").Value = 1439.5
$ws.Range("J2").Value = 1351.0555
$ws.Range("K2").Value = 1439.5
$ws.Range("L2").Value = 1351.0555
$ws.Range("M2").Value = -1326.5
$ws.Range("N2").Value = -1577.0555

$ws.Range("H61").Value = 1427.2
$ws.Range("I61").Value = 1173.4
$ws.Range("J61").Value = 2442.4
$ws.Range("K61").Value = 1173.4
$ws.Range("L61").Value = 2442.4
$ws.Range("M61").Value = -961.4000000000001
$ws.Range("N61").Value = -2866.4

$ws.Range("H116").Value = 1397.6052
$ws.Range("I116").Value = 1439.5
$ws.Range("J116").Value = 1351.0555
$ws.Range("K116").Value = 1439.5
$ws.Range("L116").Value = 1351.0555
$ws.Range("M116").Value = 854.5
$ws.Range("N116").Value = -5939.0555

$ws.Range("H135").Value = 26942.9
$ws.Range("J135").Value = 26942.9
$ws.Range("L135").Value = 26942.9
$ws.Range("N135").Value = -37082.9

$ws.Range("H136").Value = 1427.2
$ws.Range("I136").Value = 1173.4
$ws.Range("J136").Value = 2442.4
$ws.Range("K136").Value = 3520.2
$ws.Range("L136").Value = 7327.200000000001
$ws.Range("M136").Value = -970.2000000000003
$ws.Range("N136").Value = -12427.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1397.6052
$ws.Range("I3").Value = 1439.5
$ws.Range("J3").Value = 1351.0555
$ws.Range("K3").Value = 1439.5
$ws.Range("L3").Value = 1351.0555
$ws.Range("M3").Value = -1325.5
$ws.Range("N3").Value = -1579.0555

$ws.Range("H141").Value = 98616.5
$ws.Range("J141").Value = 98616.5
$ws.Range("L141").Value = 98616.5
$ws.Range("N141").Value = -108976.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2782.5957
$ws.Range("I31").Value = 1191.6842
$ws.Range("J31").Value = 3862.1428
$ws.Range("K31").Value = 1191.6842
$ws.Range("L31").Value = 3862.1428
$ws.Range("M31").Value = -896.6841999999999
$ws.Range("N31").Value = -4452.1428

$ws.Range("H34").Value = 2782.5957
$ws.Range("I34").Value = 1191.6842
$ws.Range("J34").Value = 3862.1428
$ws.Range("K34").Value = 1191.6842
$ws.Range("L34").Value = 3862.1428
$ws.Range("M34").Value = -989.6841999999999
$ws.Range("N34").Value = -4266.1428

$ws.Range("H62").Value = 2683.3333
$ws.Range("I62").Value = 2600
$ws.Range("J62").Value = 2725
$ws.Range("K62").Value = 2600
$ws.Range("L62").Value = 2725
$ws.Range("M62").Value = -1976
$ws.Range("N62").Value = -3973

$ws.Range("H65").Value = 2683.3333
$ws.Range("I65").Value = 2600
$ws.Range("J65").Value = 2725
$ws.Range("K65").Value = 13000
$ws.Range("L65").Value = 13625
$ws.Range("M65").Value = -9880
$ws.Range("N65").Value = -19865

$ws.Range("H86").Value = 6644.1304
$ws.Range("I86").Value = 3237
$ws.Range("K86").Value = 3237
$ws.Range("M86").Value = -2114

$ws.Range("H89").Value = 6644.1304
$ws.Range("I89").Value = 3237
$ws.Range("K89").Value = 16185
$ws.Range("M89").Value = -10569

$ws.Range("H99").Value = 2868
$ws.Range("I99").Value = 2837.5
$ws.Range("K99").Value = 2837.5
$ws.Range("M99").Value = -1339.5

$ws.Range("H126").Value = 2868
$ws.Range("I126").Value = 2837.5
$ws.Range("K126").Value = 8512.5
$ws.Range("M126").Value = -6042.5

$ws.Range("H141").Value = 92065.37
$ws.Range("J141").Value = 92065.37
$ws.Range("L141").Value = 92065.37
$ws.Range("N141").Value = -102425.37

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2318.2856
$ws.Range("I81").Value = 600
$ws.Range("J81").Value = 3005.6
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 9016.799999999999
$ws.Range("M81").Value = -677
$ws.Range("N81").Value = -11262.8

$ws.Range("H84").Value = 2318.2856
$ws.Range("I84").Value = 600
$ws.Range("J84").Value = 3005.6
$ws.Range("K84").Value = 5400
$ws.Range("L84").Value = 27050.4
$ws.Range("M84").Value = 216
$ws.Range("N84").Value = -38282.39999999999

$ws.Range("H113").Value = 1010508.7
$ws.Range("I113").Value = 1515580.2
$ws.Range("J113").Value = 365.4
$ws.Range("K113").Value = 4546740.6
$ws.Range("L113").Value = 1096.2
$ws.Range("M113").Value = -4544570.6
$ws.Range("N113").Value = -5436.2

$ws.Range("H131").Value = 983.36365
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 983.36365
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2950.09095
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -13030.09095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 12750
$ws.Range("J26").Value = 12750
$ws.Range("L26").Value = 12750
$ws.Range("N26").Value = -13310

$ws.Range("H50").Value = 12750
$ws.Range("J50").Value = 12750
$ws.Range("L50").Value = 12750
$ws.Range("N50").Value = -13746

$ws.Range("H80").Value = 2287.353
$ws.Range("I80").Value = 2229.9
$ws.Range("J80").Value = 2369.4285
$ws.Range("K80").Value = 2229.9
$ws.Range("L80").Value = 2369.4285
$ws.Range("M80").Value = -1231.9
$ws.Range("N80").Value = -4365.4285

$ws.Range("H83").Value = 2287.353
$ws.Range("I83").Value = 2229.9
$ws.Range("J83").Value = 2369.4285
$ws.Range("K83").Value = 11149.5
$ws.Range("L83").Value = 11847.1425
$ws.Range("M83").Value = -6157.5
$ws.Range("N83").Value = -21831.1425

$ws.Range("H102").Value = 5051654.5
$ws.Range("I102").Value = 7408371.5
$ws.Range("J102").Value = 1546.8572
$ws.Range("K102").Value = 7408371.5
$ws.Range("L102").Value = 1546.8572
$ws.Range("M102").Value = -7406749.5
$ws.Range("N102").Value = -4790.8572

$ws.Range("H126").Value = 33334368
$ws.Range("I126").Value = 55556530
$ws.Range("J126").Value = 1127.5
$ws.Range("K126").Value = 166669590
$ws.Range("L126").Value = 3382.5
$ws.Range("M126").Value = -166667120
$ws.Range("N126").Value = -8322.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3967.3333
$ws.Range("J7").Value = 5000
$ws.Range("L7").Value = 5000
$ws.Range("N7").Value = -5224

$ws.Range("H40").Value = 3102.75
$ws.Range("I40").Value = 2972.6155
$ws.Range("K40").Value = 2972.6155
$ws.Range("M40").Value = -2836.6155

$ws.Range("H61").Value = 1696.5
$ws.Range("I61").Value = 1341.25
$ws.Range("J61").Value = 3117.5
$ws.Range("K61").Value = 1341.25
$ws.Range("L61").Value = 3117.5
$ws.Range("M61").Value = -1139.25
$ws.Range("N61").Value = -3521.5

$ws.Range("H113").Value = 1696.5
$ws.Range("I113").Value = 1341.25
$ws.Range("J113").Value = 3117.5
$ws.Range("K113").Value = 1341.25
$ws.Range("L113").Value = 3117.5
$ws.Range("M113").Value = 828.75
$ws.Range("N113").Value = -7457.5

$ws.Range("H126").Value = 3967.3333
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws.Range("H136").Value = 7600.1055
$ws.Range("I136").Value = 1780.7
$ws.Range("J136").Value = 14066.111
$ws.Range("K136").Value = 5342.1
$ws.Range("L136").Value = 42198.333
$ws.Range("M136").Value = -2792.1
$ws.Range("N136").Value = -47298.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 823.5
$ws.Range("I126").Value = 823.5
$ws.Range("K126").Value = 2470.5
$ws.Range("M126").Value = -0.5

$ws.Range("H136").Value = 3109.1428
$ws.Range("I136").Value = 505.42426
$ws.Range("J136").Value = 12656.111
$ws.Range("K136").Value = 1516.27278
$ws.Range("L136").Value = 37968.333
$ws.Range("M136").Value = 1033.72722
$ws.Range("N136").Value = -43068.333
